$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.04970492700766727
$ws.Cells.Item(2, 4).Value = 0.1910414548918595
$ws.Cells.Item(2, 5).Value = 0.2625168408967085
$ws.Cells.Item(2, 6).Value = 0.6403102258717333
$ws.Cells.Item(2, 7).Value = 0.3334511935832225
$ws.Cells.Item(2, 8).Value = 0.4946654355909743
$ws.Cells.Item(2, 10).Value = 0.4453888725311117
$ws.Cells.Item(2, 13).Value = 10.85855465781913
$ws.Cells.Item(2, 15).Value = 1.589495952055742

$ws.Cells.Item(3, 3).Value = 0.0440951706950159
$ws.Cells.Item(3, 4).Value = 0.1902076448927517
$ws.Cells.Item(3, 5).Value = 0.2439416287165628
$ws.Cells.Item(3, 6).Value = 0.671216691506082
$ws.Cells.Item(3, 7).Value = 0.3420505621732701
$ws.Cells.Item(3, 8).Value = 0.5074969017347186
$ws.Cells.Item(3, 10).Value = 0.401270692052293
$ws.Cells.Item(3, 13).Value = 9.524621878947414
$ws.Cells.Item(3, 15).Value = 1.634078423046745

$ws.Cells.Item(4, 3).Value = 0.04066862791573556
$ws.Cells.Item(4, 4).Value = 0.1898913335161438
$ws.Cells.Item(4, 5).Value = 0.2327722267321448
$ws.Cells.Item(4, 6).Value = 0.6915747250045996
$ws.Cells.Item(4, 7).Value = 0.3481778328296841
$ws.Cells.Item(4, 8).Value = 0.5160105957491083
$ws.Cells.Item(4, 10).Value = 0.3744832907197804
$ws.Cells.Item(4, 13).Value = 8.702353442672461
$ws.Cells.Item(4, 15).Value = 1.664563336063097

$ws.Cells.Item(5, 3).Value = 0.03927665258281365
$ws.Cells.Item(5, 4).Value = 0.1898110255528707
$ws.Cells.Item(5, 5).Value = 0.2282781085516419
$ws.Cells.Item(5, 6).Value = 0.7002122721071267
$ws.Cells.Item(5, 7).Value = 0.3508834230483089
$ws.Cells.Item(5, 8).Value = 0.5196378245934028
$ws.Cells.Item(5, 10).Value = 0.3636394778143881
$ws.Cells.Item(5, 13).Value = 8.366434835913765
$ws.Cells.Item(5, 15).Value = 1.677755304737744

$ws.Cells.Item(6, 3).Value = 0.03904577642511242
$ws.Cells.Item(6, 4).Value = 0.1898006082967783
$ws.Cells.Item(6, 5).Value = 0.2275352865772646
$ws.Cells.Item(6, 6).Value = 0.7016669843573027
$ws.Cells.Item(6, 7).Value = 0.3513451608975231
$ws.Cells.Item(6, 8).Value = 0.5202496009521269
$ws.Cells.Item(6, 10).Value = 0.3618431386420298
$ws.Cells.Item(6, 13).Value = 8.31060455802907
$ws.Cells.Item(6, 15).Value = 1.679991889934726

$ws.Cells.Item(7, 3).Value = 0.04064983772161668
$ws.Cells.Item(7, 4).Value = 0.1898900544818076
$ws.Cells.Item(7, 5).Value = 0.2327113869375381
$ws.Cells.Item(7, 6).Value = 0.6916898390391566
$ws.Cells.Item(7, 7).Value = 0.3482134821743514
$ws.Cells.Item(7, 8).Value = 0.5160588772500176
$ws.Cells.Item(7, 10).Value = 0.374336759081018
$ws.Cells.Item(7, 13).Value = 8.697826548810383
$ws.Cells.Item(7, 15).Value = 1.664738150777708

$ws.Cells.Item(8, 3).Value = 0.04776689596079109
$ws.Cells.Item(8, 4).Value = 0.1907129715356604
$ws.Cells.Item(8, 5).Value = 0.2560620236514666
$ws.Cells.Item(8, 6).Value = 0.6506765388418145
$ws.Cells.Item(8, 7).Value = 0.3362377791251987
$ws.Cells.Item(8, 8).Value = 0.498956774987434
$ws.Cells.Item(8, 10).Value = 0.4301121756318764
$ws.Cells.Item(8, 13).Value = 10.39926736874446
$ws.Cells.Item(8, 15).Value = 1.604214319534449

$ws.Cells.Item(9, 3).Value = 0.06187173935275325
$ws.Cells.Item(9, 4).Value = 0.193908903109886
$ws.Cells.Item(9, 5).Value = 0.3038145499261447
$ws.Cells.Item(9, 6).Value = 0.581473521005023
$ws.Cells.Item(9, 7).Value = 0.3196730856583656
$ws.Cells.Item(9, 8).Value = 0.4705430926911305
$ws.Cells.Item(9, 10).Value = 0.5420557495075684
$ws.Cells.Item(9, 13).Value = 13.71173061222999
$ws.Cells.Item(9, 15).Value = 1.510807200970476

$ws.Cells.Item(10, 3).Value = 0.07233530378638875
$ws.Cells.Item(10, 4).Value = 0.1972664231525272
$ws.Cells.Item(10, 5).Value = 0.3402336975439084
$ws.Cells.Item(10, 6).Value = 0.5378369645073207
$ws.Cells.Item(10, 7).Value = 0.3119957722768021
$ws.Cells.Item(10, 8).Value = 0.4529068009134249
$ws.Cells.Item(10, 10).Value = 0.6261400426023158
$ws.Cells.Item(10, 13).Value = 16.13346934510241
$ws.Cells.Item(10, 15).Value = 1.45841724915573

$ws.Cells.Item(11, 3).Value = 0.07711970544556834
$ws.Cells.Item(11, 4).Value = 0.1990237830018344
$ws.Cells.Item(11, 5).Value = 0.357125257179689
$ws.Cells.Item(11, 6).Value = 0.5196306720905071
$ws.Cells.Item(11, 7).Value = 0.3095390793906176
$ws.Cells.Item(11, 8).Value = 0.4456123025196774
$ws.Cells.Item(11, 10).Value = 0.6648577257375621
$ws.Cells.Item(11, 13).Value = 17.23330324063608
$ws.Cells.Item(11, 15).Value = 1.438290422727675

$ws.Cells.Item(12, 3).Value = 0.07893513201051405
$ws.Cells.Item(12, 4).Value = 0.1997232696770084
$ws.Cells.Item(12, 5).Value = 0.363571208628727
$ws.Cells.Item(12, 6).Value = 0.5129796881554398
$ws.Cells.Item(12, 7).Value = 0.3087628157352071
$ws.Cells.Item(12, 8).Value = 0.4429569477018163
$ws.Cells.Item(12, 10).Value = 0.6795920460965021
$ws.Cells.Item(12, 13).Value = 17.64957992839527
$ws.Cells.Item(12, 15).Value = 1.431217064765548

$ws.Cells.Item(13, 3).Value = 0.07854398127662421
$ws.Cells.Item(13, 4).Value = 0.199571094414992
$ws.Cells.Item(13, 5).Value = 0.3621807093968954
$ws.Cells.Item(13, 6).Value = 0.5144011621493334
$ws.Cells.Item(13, 7).Value = 0.3089230651745254
$ws.Cells.Item(13, 8).Value = 0.4435240357576191
$ws.Cells.Item(13, 10).Value = 0.6764154170933239
$ws.Cells.Item(13, 13).Value = 17.55993552377959
$ws.Cells.Item(13, 15).Value = 1.432715808146298

$ws.Cells.Item(14, 3).Value = 0.07726898705227825
$ws.Cells.Item(14, 4).Value = 0.199080642947763
$ws.Cells.Item(14, 5).Value = 0.3576545607419916
$ws.Cells.Item(14, 6).Value = 0.5190785745678497
$ws.Cells.Item(14, 7).Value = 0.3094720991386311
$ws.Cells.Item(14, 8).Value = 0.4453916895620011
$ws.Cells.Item(14, 10).Value = 0.6660684353540205
$ws.Cells.Item(14, 13).Value = 17.26755437063053
$ws.Cells.Item(14, 15).Value = 1.437697416790343

$ws.Cells.Item(15, 3).Value = 0.07648849999155516
$ws.Cells.Item(15, 4).Value = 0.1987846855571433
$ws.Cells.Item(15, 5).Value = 0.3548886961670661
$ws.Cells.Item(15, 6).Value = 0.5219755228767369
$ws.Cells.Item(15, 7).Value = 0.3098286107656776
$ws.Cells.Item(15, 8).Value = 0.4465496700508851
$ws.Cells.Item(15, 10).Value = 0.6597402595491815
$ws.Cells.Item(15, 13).Value = 17.08843734296795
$ws.Cells.Item(15, 15).Value = 1.440820653459184

$ws.Cells.Item(16, 3).Value = 0.07202313774773472
$ws.Cells.Item(16, 4).Value = 0.1971562827127968
$ws.Cells.Item(16, 5).Value = 0.3391365595189058
$ws.Cells.Item(16, 6).Value = 0.5390604613847785
$ws.Cells.Item(16, 7).Value = 0.3121775979767705
$ws.Cells.Item(16, 8).Value = 0.4533983549988108
$ws.Cells.Item(16, 10).Value = 0.6236196224449202
$ws.Cells.Item(16, 13).Value = 16.06156051752913
$ws.Cells.Item(16, 15).Value = 1.459808458415608

$ws.Cells.Item(17, 3).Value = 0.06929016163610413
$ws.Cells.Item(17, 4).Value = 0.1962168694057027
$ws.Cells.Item(17, 5).Value = 0.3295581641183958
$ws.Cells.Item(17, 6).Value = 0.5499675073503525
$ws.Cells.Item(17, 7).Value = 0.3138874989245579
$ws.Cells.Item(17, 8).Value = 0.4577878851596893
$ws.Cells.Item(17, 10).Value = 0.6015843195198158
$ws.Cells.Item(17, 13).Value = 15.43117500997715
$ws.Cells.Item(17, 15).Value = 1.472416805215147

$ws.Cells.Item(18, 3).Value = 0.06772051982031257
$ws.Cells.Item(18, 4).Value = 0.1956981334085413
$ws.Cells.Item(18, 5).Value = 0.3240792676763533
$ws.Cells.Item(18, 6).Value = 0.5563952895404185
$ws.Cells.Item(18, 7).Value = 0.3149681685354437
$ws.Cells.Item(18, 8).Value = 0.4603810490670668
$ws.Cells.Item(18, 10).Value = 0.5889538256062394
$ws.Cells.Item(18, 13).Value = 15.06841823235516
$ws.Cells.Item(18, 15).Value = 1.4800166476698

$ws.Cells.Item(19, 3).Value = 0.06718945530398912
$ws.Cells.Item(19, 4).Value = 0.1955261787459932
$ws.Cells.Item(19, 5).Value = 0.3222293352410759
$ws.Cells.Item(19, 6).Value = 0.5585979255983737
$ws.Cells.Item(19, 7).Value = 0.3153506018991834
$ws.Cells.Item(19, 8).Value = 0.4612707356102703
$ws.Cells.Item(19, 10).Value = 0.5846846871148159
$ws.Cells.Item(19, 13).Value = 14.94556322510334
$ws.Cells.Item(19, 15).Value = 1.48264910458019

$ws.Cells.Item(20, 3).Value = 0.06958085308663442
$ws.Cells.Item(20, 4).Value = 0.1963146304990602
$ws.Cells.Item(20, 5).Value = 0.3305746402137117
$ws.Cells.Item(20, 6).Value = 0.5487904107168973
$ws.Cells.Item(20, 7).Value = 0.3136953827463458
$ws.Cells.Item(20, 8).Value = 0.4573135150753842
$ws.Cells.Item(20, 10).Value = 0.603925460468389
$ws.Cells.Item(20, 13).Value = 15.49829841972877
$ws.Cells.Item(20, 15).Value = 1.471038513129713

$ws.Cells.Item(21, 3).Value = 0.07764338266547099
$ws.Cells.Item(21, 4).Value = 0.1992237695660179
$ws.Cells.Item(21, 5).Value = 0.3589826336613413
$ws.Cells.Item(21, 6).Value = 0.5176980450626019
$ws.Cells.Item(21, 7).Value = 0.3093066130571316
$ws.Cells.Item(21, 8).Value = 0.444840195205515
$ws.Cells.Item(21, 10).Value = 0.6691055733186033
$ws.Cells.Item(21, 13).Value = 17.35343893119938
$ws.Cells.Item(21, 15).Value = 1.436219194659685

$ws.Cells.Item(22, 3).Value = 0.08293422479682988
$ws.Cells.Item(22, 4).Value = 0.2013237899492424
$ws.Cells.Item(22, 5).Value = 0.3778387928358171
$ws.Cells.Item(22, 6).Value = 0.4987990956293089
$ws.Cells.Item(22, 7).Value = 0.3073385566041651
$ws.Cells.Item(22, 8).Value = 0.4373123724451915
$ws.Cells.Item(22, 10).Value = 0.7121315887567903
$ws.Cells.Item(22, 13).Value = 18.56470690915631
$ws.Cells.Item(22, 15).Value = 1.416665638479998

$ws.Cells.Item(23, 3).Value = 0.08010838049614222
$ws.Cells.Item(23, 4).Value = 0.2001844602789191
$ws.Cells.Item(23, 5).Value = 0.3677473889161718
$ws.Cells.Item(23, 6).Value = 0.50875341193629
$ws.Cells.Item(23, 7).Value = 0.3083048619992184
$ws.Cells.Item(23, 8).Value = 0.4412722673588974
$ws.Cells.Item(23, 10).Value = 0.6891267784953357
$ws.Cells.Item(23, 13).Value = 17.91831767668367
$ws.Cells.Item(23, 15).Value = 1.426803518025707

$ws.Cells.Item(24, 3).Value = 0.06944942656460285
$ws.Cells.Item(24, 4).Value = 0.1962703662934047
$ws.Cells.Item(24, 5).Value = 0.330115004835605
$ws.Cells.Item(24, 6).Value = 0.5493220870575826
$ws.Cells.Item(24, 7).Value = 0.313781934775335
$ws.Cells.Item(24, 8).Value = 0.4575277612771771
$ws.Cells.Item(24, 10).Value = 0.6028669130271851
$ws.Cells.Item(24, 13).Value = 15.46795298860866
$ws.Cells.Item(24, 15).Value = 1.471660546431309

$ws.Cells.Item(25, 3).Value = 0.05803908011424141
$ws.Cells.Item(25, 4).Value = 0.1928704803881658
$ws.Cells.Item(25, 5).Value = 0.290672409937514
$ws.Cells.Item(25, 6).Value = 0.5989558930814383
$ws.Cells.Item(25, 7).Value = 0.3233853968855414
$ws.Cells.Item(25, 8).Value = 0.4776693665833989
$ws.Cells.Item(25, 10).Value = 0.511470065248119
$ws.Cells.Item(25, 13).Value = 12.81793408937523
$ws.Cells.Item(25, 15).Value = 1.533285433783391
